# Update Name of Algo
# Applies corrected KNN-imputed values to columns B and C for the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C7").Value = -13.566
$ws.Range("B8").Value = 6.779999999999999
$ws.Range("B10").Value = 6.536
$ws.Range("B12").Value = 5.166
$ws.Range("C14").Value = -12.096
$ws.Range("C15").Value = -13.636
$ws.Range("B18").Value = 5.761999999999999
$ws.Range("C18").Value = -12.379
$ws.Range("C20").Value = -12.252
$ws.Range("B25").Value = 6.985000000000001
$ws.Range("C29").Value = -12.423
$ws.Range("C30").Value = -12.952
$ws.Range("C31").Value = -12.943
$ws.Range("C35").Value = -11.986
$ws.Range("B37").Value = 7.997
$ws.Range("C40").Value = -12.782
$ws.Range("C44").Value = -12.133
$ws.Range("C50").Value = -13.326
$ws.Range("C54").Value = -12.72
$ws.Range("B55").Value = 5.194000000000001
$ws.Range("B68").Value = 5.431
$ws.Range("C68").Value = -11.177
$ws.Range("C76").Value = -13.366
$ws.Range("B77").Value = 5.737
$ws.Range("B78").Value = 7.696
$ws.Range("B79").Value = 5.353
$ws.Range("B80").Value = 8.087999999999999
$ws.Range("B81").Value = 6.004
$ws.Range("B82").Value = 6.009
$ws.Range("B84").Value = 6.207000000000001
$ws.Range("C87").Value = -12.826
$ws.Range("C88").Value = -13.035
$ws.Range("C92").Value = -11.301
$ws.Range("C96").Value = -12.934
$ws.Range("C98").Value = -13.23
$ws.Range("B101").Value = 9.032
$ws.Range("C101").Value = -12.442
$ws.Range("B102").Value = 7.286
$ws.Range("C102").Value = -12.981
